$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9691714836223507
$ws.Range("C2").Value = 0.8613013698630136
$ws.Range("D2").Value = 0.9120580235720761

$ws.Range("B3").Value = 0.3467741935483871
$ws.Range("C3").Value = 0.7288135593220338
$ws.Range("D3").Value = 0.4699453551912569

$ws.Range("B4").Value = 0.8491446345256609
$ws.Range("C4").Value = 0.8491446345256609
$ws.Range("D4").Value = 0.8491446345256609
$ws.Range("E4").Value = 0.8491446345256609

$ws.Range("B5").Value = 0.6579728385853689
$ws.Range("C5").Value = 0.7950574645925237
$ws.Range("D5").Value = 0.6910016893816665

$ws.Range("B6").Value = 0.9120619344553774
$ws.Range("C6").Value = 0.8491446345256609
$ws.Range("D6").Value = 0.8714909202525297
